# "Stat - added last stat"
#
# TimeTracker sheet: each data row is one day (A=date-of-month, D=day label,
# E=daily total formula, F:I=per-entry hours). This edit:
#   1. Adds one more hour entry (G10=1) to the existing "Utorok 30.4" row,
#      bumping its total (E10) from 2 -> 3.
#   2. Adds a brand-new day row 11 for "Streda 1.5" with its own running
#      total formula (starts at 0, no entries logged yet).
#   3. The week grand total in E13 (=SUM(E2:E12)) picks up the new values
#      automatically (44.5 -> 45.5).
#   4. Leaves the selection on E11, matching where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Tuesday 30.4 (row 10) gets one more logged hour.
$ws.Range("G10").Value = 1

# 2) New day: Wednesday 1.5 (row 11).
$ws.Range("D11").Value = "Streda 1.5"
$ws.Range("E11").Formula = "=SUM(F11:I11)"

# Bold the new/updated daily-total cells to match the style already used
# on the other "total" cells in column E (E7:E9).
$ws.Range("E10").Font.Bold = $true
$ws.Range("E11").Font.Bold = $true

# 3) E13 already holds =SUM(E2:E12), so it recalculates to 45.5 on its own.

# 4) Leave the cursor on the new total cell.
$ws.Range("E11").Select()
